$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("week1")

# Set the actual time length to complete for "DQ1 response 2" (row 8, column C)
# Value represents 10 minutes expressed as a fraction of a day ([h]:mm format)
$ws.Range("C8").Value = 0.006944444444444444

# Recalculate so the Total formula (C20) picks up the new value
$excel.Calculate()

# Move the active selection to C9, as in the saved workbook
$ws.Range("C9").Select()
